# GDPT rmse-z by c-min supplementary figure data update
# (added a script to compare IDPT, SPCT and GDPT)
#
# Changes applied:
#  1. Column G (avg_num_meas_by_z) is recomputed as num_meas / 63 for every
#     data row (rows 2-50), replacing the old (much smaller) values.
#  2. Column D (rmse_z_ga) is updated for rows 22-50 with refreshed GDPT
#     RMSE-z values from the re-run comparison.
#  3. Row 51 (cmin = 1, the empty tail bin) no longer has a numeric D51
#     value - it becomes blank, matching the already-blank C51 cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. avg_num_meas_by_z = num_meas / 63 for rows 2 through 50 ---
for ($r = 2; $r -le 50; $r++) {
    $numMeas = $ws.Cells.Item($r, 5).Value2
    $ws.Cells.Item($r, 7).Value = $numMeas / 63
}

# --- 2. refreshed rmse_z_ga values (column D) for rows 22-50 ---
$dValues = @{
    22 = 2.136793603657316
    23 = 2.136280133796781
    24 = 2.13555450967014
    25 = 2.127589679888585
    26 = 2.127329707785147
    27 = 2.124988642338044
    28 = 2.126195651048691
    29 = 2.124861618068633
    30 = 2.121556542345739
    31 = 2.119864029668172
    32 = 2.119958663608569
    33 = 2.113222016897527
    34 = 2.111660257409673
    35 = 2.112066289773897
    36 = 2.105810207186438
    37 = 2.102902641161698
    38 = 2.104468209883738
    39 = 2.099912607935125
    40 = 2.091878885839314
    41 = 2.079396135689238
    42 = 2.065359644865959
    43 = 2.042622446870265
    44 = 2.034839747563011
    45 = 2.021975247708756
    46 = 1.997742026431591
    47 = 1.95693254851187
    48 = 1.835499263639429
    49 = 1.504352486364895
    50 = 1.270685148540861
}
foreach ($r in $dValues.Keys) {
    $ws.Cells.Item($r, 4).Value = $dValues[$r]
}

# --- 3. D51 becomes blank (like C51) instead of holding the old constant ---
# A leading single-quote forces an empty text entry (matching C51's empty
# inline-string cell) rather than simply clearing the cell to a blank
# numeric cell.
$ws.Range("D51").Value = "'"
